$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1300
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1300
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1300
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1952
$ws.Range("H132").Value = 2382.5715
$ws.Range("I132").Value = 2382.5715
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7147.7145
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4617.7145
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 1346.875
$ws.Range("I137").Value = 1328
$ws.Range("K137").Value = 3984
$ws.Range("M137").Value = -1434
$ws.Range("H138").Value = 197
$ws.Range("I138").Value = 197
$ws.Range("K138").Value = 591
$ws.Range("M138").Value = 4549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2301
$ws.Range("I2").Value = 952
$ws.Range("K2").Value = 952
$ws.Range("M2").Value = -839
$ws.Range("H3").Value = 1006
$ws.Range("J3").Value = 1006
$ws.Range("L3").Value = 1006
$ws.Range("N3").Value = -1236
$ws.Range("H32").Value = 4548.35
$ws.Range("I32").Value = 3599.9412
$ws.Range("J32").Value = 9922.666999999999
$ws.Range("K32").Value = 3599.9412
$ws.Range("L32").Value = 9922.666999999999
$ws.Range("M32").Value = -3312.9412
$ws.Range("N32").Value = -10496.667
$ws.Range("H61").Value = 1411.9333
$ws.Range("I61").Value = 1411.9333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1411.9333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1199.9333
$ws.Range("N61").ClearContents()
$ws.Range("H116").Value = 2301
$ws.Range("I116").Value = 952
$ws.Range("K116").Value = 952
$ws.Range("M116").Value = 1342
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H136").Value = 1411.9333
$ws.Range("I136").Value = 1411.9333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4235.7999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1685.7999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2301
$ws.Range("I3").Value = 952
$ws.Range("K3").Value = 952
$ws.Range("M3").Value = -838
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H86").Value = 4448.1665
$ws.Range("I86").Value = 3197.25
$ws.Range("J86").Value = 6950
$ws.Range("K86").Value = 3197.25
$ws.Range("L86").Value = 6950
$ws.Range("M86").Value = -2074.25
$ws.Range("N86").Value = -9196
$ws.Range("H89").Value = 4448.1665
$ws.Range("I89").Value = 3197.25
$ws.Range("J89").Value = 6950
$ws.Range("K89").Value = 15986.25
$ws.Range("L89").Value = 34750
$ws.Range("M89").Value = -10370.25
$ws.Range("N89").Value = -45982
$ws.Range("H105").Value = 4190.3335
$ws.Range("I105").Value = 3776.625
$ws.Range("K105").Value = 3776.625
$ws.Range("M105").Value = -2029.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H115").Value = 56000
$ws.Range("J115").Value = 56000
$ws.Range("L115").Value = 56000
$ws.Range("N115").Value = -58350
$ws.Range("H134").Value = 2575.3333
$ws.Range("I134").Value = 2575.3333
$ws.Range("K134").Value = 7725.999899999999
$ws.Range("M134").Value = -5190.999899999999
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 111222580
$ws.Range("I4").Value = 457.57144
$ws.Range("K4").Value = 1372.71432
$ws.Range("M4").Value = -1260.71432
$ws.Range("H9").Value = 3000
$ws.Range("J9").Value = 3000
$ws.Range("L9").Value = 9000
$ws.Range("N9").Value = -9448
$ws.Range("H10").Value = 63.666668
$ws.Range("I10").Value = 20.5
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 61.5
$ws.Range("L10").Value = 450
$ws.Range("M10").Value = 77.5
$ws.Range("N10").Value = -728
$ws.Range("H15").Value = 159.14285
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 177.33333
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 531.99999
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -811.99999
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -1027
$ws.Range("H39").Value = 3200
$ws.Range("J39").Value = 3200
$ws.Range("L39").Value = 9600
$ws.Range("N39").Value = -10188
$ws.Range("H112").Value = 42749.25
$ws.Range("I112").Value = 21000
$ws.Range("J112").Value = 49999
$ws.Range("K112").Value = 63000
$ws.Range("L112").Value = 149997
$ws.Range("M112").Value = -61892
$ws.Range("N112").Value = -152213
$ws.Range("H131").Value = 1182.8572
$ws.Range("J131").Value = 1525
$ws.Range("L131").Value = 4575
$ws.Range("N131").Value = -14655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 13004.5
$ws.Range("J40").Value = 13004.5
$ws.Range("L40").Value = 13004.5
$ws.Range("N40").Value = -13306.5
$ws.Range("H102").Value = 1688.9231
$ws.Range("I102").Value = 1666.4
$ws.Range("K102").Value = 1666.4
$ws.Range("M102").Value = -44.40000000000009
$ws.Range("H136").Value = 25000
$ws.Range("J136").Value = 35000
$ws.Range("L136").Value = 105000
$ws.Range("N136").Value = -110100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H40").Value = 1966.6666
$ws.Range("I40").Value = 1700
$ws.Range("K40").Value = 1700
$ws.Range("M40").Value = -1564
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H110").Value = 50214.668
$ws.Range("J110").Value = 50214.668
$ws.Range("L110").Value = 50214.668
$ws.Range("N110").Value = -58394.668
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 57857.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 57857.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 57857.5
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -68057.5
